$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cell F1 with same style as E1 (bold/border/center/top)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate time_taken values for data rows
$ws.Range("F2").Value = "2021-10-05 10:50:20.396835"
$ws.Range("F3").Value = "2021-10-05 10:50:20.396847"
$ws.Range("F4").Value = "2021-10-05 10:50:20.396850"
$ws.Range("F5").Value = "2021-10-05 10:50:20.396853"
$ws.Range("F6").Value = "2021-10-05 10:50:20.396856"
$ws.Range("F7").Value = "2021-10-05 10:50:20.396859"
$ws.Range("F8").Value = "2021-10-05 10:50:20.396861"
$ws.Range("F9").Value = "2021-10-05 10:50:20.396864"
$ws.Range("F10").Value = "2021-10-05 10:50:20.396866"
$ws.Range("F11").Value = "2021-10-05 10:50:20.396869"
$ws.Range("F12").Value = "2021-10-05 10:50:20.396871"
$ws.Range("F13").Value = "2021-10-05 10:50:20.396874"
$ws.Range("F14").Value = "2021-10-05 10:50:20.396876"
$ws.Range("F15").Value = "2021-10-05 10:50:20.396879"
$ws.Range("F16").Value = "2021-10-05 10:50:20.396882"
$ws.Range("F17").Value = "2021-10-05 10:50:20.396884"
$ws.Range("F18").Value = "2021-10-05 10:50:20.396887"
$ws.Range("F19").Value = "2021-10-05 10:50:20.396890"
$ws.Range("F20").Value = "2021-10-05 10:50:20.396893"
$ws.Range("F21").Value = "2021-10-05 10:50:20.396895"
$ws.Range("F22").Value = "2021-10-05 10:50:20.396898"
$ws.Range("F23").Value = "2021-10-05 10:50:20.396900"
$ws.Range("F24").Value = "2021-10-05 10:50:20.396903"
$ws.Range("F25").Value = "2021-10-05 10:50:20.396905"
$ws.Range("F26").Value = "2021-10-05 10:50:20.396908"
$ws.Range("F27").Value = "2021-10-05 10:50:20.396911"
$ws.Range("F28").Value = "2021-10-05 10:50:20.396913"
$ws.Range("F29").Value = "2021-10-05 10:50:20.396916"
$ws.Range("F30").Value = "2021-10-05 10:50:20.396918"
$ws.Range("F31").Value = "2021-10-05 10:50:20.396921"
$ws.Range("F32").Value = "2021-10-05 10:50:20.396923"
$ws.Range("F33").Value = "2021-10-05 10:50:20.396926"
$ws.Range("F34").Value = "2021-10-05 10:50:20.396929"
$ws.Range("F35").Value = "2021-10-05 10:50:20.396931"
$ws.Range("F36").Value = "2021-10-05 10:50:20.396934"
$ws.Range("F37").Value = "2021-10-05 10:50:20.396936"
$ws.Range("F38").Value = "2021-10-05 10:50:20.396939"
$ws.Range("F39").Value = "2021-10-05 10:50:20.396942"
$ws.Range("F40").Value = "2021-10-05 10:50:20.396944"
$ws.Range("F41").Value = "2021-10-05 10:50:20.396947"
$ws.Range("F42").Value = "2021-10-05 10:50:20.396950"
$ws.Range("F43").Value = "2021-10-05 10:50:20.396952"
$ws.Range("F44").Value = "2021-10-05 10:50:20.396955"
$ws.Range("F45").Value = "2021-10-05 10:50:20.396958"
$ws.Range("F46").Value = "2021-10-05 10:50:20.396960"
$ws.Range("F47").Value = "2021-10-05 10:50:20.396963"
$ws.Range("F48").Value = "2021-10-05 10:50:20.396965"
$ws.Range("F49").Value = "2021-10-05 10:50:20.396968"
$ws.Range("F50").Value = "2021-10-05 10:50:20.396970"
$ws.Range("F51").Value = "2021-10-05 10:50:20.396973"
$ws.Range("F52").Value = "2021-10-05 10:50:20.396975"
$ws.Range("F53").Value = "2021-10-05 10:50:20.396978"
$ws.Range("F54").Value = "2021-10-05 10:50:20.396981"
$ws.Range("F55").Value = "2021-10-05 10:50:20.396983"
$ws.Range("F56").Value = "2021-10-05 10:50:20.396986"
$ws.Range("F57").Value = "2021-10-05 10:50:20.396989"
$ws.Range("F58").Value = "2021-10-05 10:50:20.396991"
$ws.Range("F59").Value = "2021-10-05 10:50:20.396994"
$ws.Range("F60").Value = "2021-10-05 10:50:20.396996"
$ws.Range("F61").Value = "2021-10-05 10:50:20.396999"
$ws.Range("F62").Value = "2021-10-05 10:50:20.397002"
$ws.Range("F63").Value = "2021-10-05 10:50:20.397004"
$ws.Range("F64").Value = "2021-10-05 10:50:20.397007"
$ws.Range("F65").Value = "2021-10-05 10:50:20.397009"
$ws.Range("F66").Value = "2021-10-05 10:50:20.397013"
$ws.Range("F67").Value = "2021-10-05 10:50:20.397016"
$ws.Range("F68").Value = "2021-10-05 10:50:20.397019"
$ws.Range("F69").Value = "2021-10-05 10:50:20.397021"
$ws.Range("F70").Value = "2021-10-05 10:50:20.397024"
$ws.Range("F71").Value = "2021-10-05 10:50:20.397026"
$ws.Range("F72").Value = "2021-10-05 10:50:20.397029"
$ws.Range("F73").Value = "2021-10-05 10:50:20.397031"
$ws.Range("F74").Value = "2021-10-05 10:50:20.397034"
$ws.Range("F75").Value = "2021-10-05 10:50:20.397036"
$ws.Range("F76").Value = "2021-10-05 10:50:20.397039"
$ws.Range("F77").Value = "2021-10-05 10:50:20.397041"
$ws.Range("F78").Value = "2021-10-05 10:50:20.397046"
$ws.Range("F79").Value = "2021-10-05 10:50:20.397049"
$ws.Range("F80").Value = "2021-10-05 10:50:20.397052"
$ws.Range("F81").Value = "2021-10-05 10:50:20.397054"
$ws.Range("F82").Value = "2021-10-05 10:50:20.397057"
$ws.Range("F83").Value = "2021-10-05 10:50:20.397059"
$ws.Range("F84").Value = "2021-10-05 10:50:20.397062"
$ws.Range("F85").Value = "2021-10-05 10:50:20.397064"
$ws.Range("F86").Value = "2021-10-05 10:50:20.397067"
$ws.Range("F87").Value = "2021-10-05 10:50:20.397070"
$ws.Range("F88").Value = "2021-10-05 10:50:20.397072"
$ws.Range("F89").Value = "2021-10-05 10:50:20.397075"
$ws.Range("F90").Value = "2021-10-05 10:50:20.397077"
$ws.Range("F91").Value = "2021-10-05 10:50:20.397080"
$ws.Range("F92").Value = "2021-10-05 10:50:20.397083"
$ws.Range("F93").Value = "2021-10-05 10:50:20.397085"
$ws.Range("F94").Value = "2021-10-05 10:50:20.397089"
$ws.Range("F95").Value = "2021-10-05 10:50:20.397092"
$ws.Range("F96").Value = "2021-10-05 10:50:20.397095"
$ws.Range("F97").Value = "2021-10-05 10:50:20.397097"
$ws.Range("F98").Value = "2021-10-05 10:50:20.397100"
$ws.Range("F99").Value = "2021-10-05 10:50:20.397103"
$ws.Range("F100").Value = "2021-10-05 10:50:20.397105"
